$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 6")
Write-Host $ws.Name
